$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1079.8889
$ws.Range("I9").Value = 453.33334
$ws.Range("K9").Value = 453.33334
$ws.Range("M9").Value = -284.33334
$ws.Range("H16").Value = 29999
$ws.Range("J16").Value = 29999
$ws.Range("L16").Value = 29999
$ws.Range("N16").Value = -30459
$ws.Range("H33").Value = 1779.6923
$ws.Range("I33").Value = 418.16666
$ws.Range("J33").Value = 2946.7144
$ws.Range("K33").Value = 418.16666
$ws.Range("L33").Value = 2946.7144
$ws.Range("M33").Value = -189.16666
$ws.Range("N33").Value = -3404.7144
$ws.Range("H40").Value = 10143.1
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10143.1
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10143.1
$ws.Range("N40").Value = -10493.1
$ws.Range("M40").ClearContents()
$ws.Range("H43").Value = 6099.926
$ws.Range("J43").Value = 3462.25
$ws.Range("L43").Value = 3462.25
$ws.Range("N43").Value = -3600.25
$ws.Range("H69").Value = 16207.477
$ws.Range("I69").Value = 8989.799999999999
$ws.Range("J69").Value = 18463
$ws.Range("K69").Value = 26969.4
$ws.Range("L69").Value = 55389
$ws.Range("M69").Value = -26095.4
$ws.Range("N69").Value = -57137
$ws.Range("H72").Value = 16207.477
$ws.Range("I72").Value = 8989.799999999999
$ws.Range("J72").Value = 18463
$ws.Range("K72").Value = 80908.2
$ws.Range("L72").Value = 166167
$ws.Range("M72").Value = -76540.2
$ws.Range("N72").Value = -174903
$ws.Range("H74").Value = 10499.5
$ws.Range("I74").Value = 3998
$ws.Range("K74").Value = 3998
$ws.Range("M74").Value = -3062
$ws.Range("H77").Value = 10499.5
$ws.Range("I77").Value = 3998
$ws.Range("K77").Value = 19990
$ws.Range("M77").Value = -15310
$ws.Range("H100").Value = 3778.1428
$ws.Range("I100").Value = 3581.3635
$ws.Range("K100").Value = 3581.3635
$ws.Range("M100").Value = -3040.3635
$ws.Range("H101").Value = 1219
$ws.Range("I101").Value = 1098
$ws.Range("J101").Value = 1299.6666
$ws.Range("K101").Value = 3294
$ws.Range("L101").Value = 3898.9998
$ws.Range("M101").Value = -1672
$ws.Range("N101").Value = -7142.9998
$ws.Range("H112").Value = 6704
$ws.Range("J112").Value = 8965.799999999999
$ws.Range("L112").Value = 26897.4
$ws.Range("N112").Value = -29113.4
$ws.Range("H116").Value = 17681.445
$ws.Range("I116").Value = 17460.23
$ws.Range("K116").Value = 17460.23
$ws.Range("M116").Value = -14018.23
$ws.Range("H117").Value = 78666.336
$ws.Range("J117").Value = 78666.336
$ws.Range("L117").Value = 78666.336
$ws.Range("N117").Value = -87844.336
$ws.Range("H132").Value = 2810.1333
$ws.Range("I132").Value = 2786.36
$ws.Range("K132").Value = 8359.08
$ws.Range("M132").Value = -5829.08
$ws.Range("H137").Value = 2270.8928
$ws.Range("I137").Value = 2204.1304
$ws.Range("K137").Value = 6612.3912
$ws.Range("M137").Value = -4062.3912
$ws.Range("H138").Value = 3906.973
$ws.Range("I138").Value = 3301.1628
$ws.Range("J138").Value = 4747.2905
$ws.Range("K138").Value = 9903.4884
$ws.Range("L138").Value = 14241.8715
$ws.Range("M138").Value = -4763.4884
$ws.Range("N138").Value = -24521.8715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 934.4286
$ws.Range("I2").Value = 963.1053000000001
$ws.Range("J2").Value = 662
$ws.Range("K2").Value = 963.1053000000001
$ws.Range("L2").Value = 662
$ws.Range("M2").Value = -850.1053000000001
$ws.Range("N2").Value = -888
$ws.Range("H32").Value = 7105.3823
$ws.Range("I32").Value = 3301.9648
$ws.Range("K32").Value = 3301.9648
$ws.Range("M32").Value = -3014.9648
$ws.Range("H45").Value = 4268.4
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4268.4
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4268.4
$ws.Range("N45").Value = -5022.4
$ws.Range("M45").ClearContents()
$ws.Range("H74").Value = 16671930
$ws.Range("I74").Value = 31252036
$ws.Range("K74").Value = 31252036
$ws.Range("M74").Value = -31251162
$ws.Range("H77").Value = 16671930
$ws.Range("I77").Value = 31252036
$ws.Range("K77").Value = 156260180
$ws.Range("M77").Value = -156255812
$ws.Range("H110").Value = 2772.5789
$ws.Range("I110").Value = 2787.7778
$ws.Range("J110").Value = 2499
$ws.Range("K110").Value = 2787.7778
$ws.Range("L110").Value = 2499
$ws.Range("M110").Value = -742.7777999999998
$ws.Range("N110").Value = -6589
$ws.Range("H116").Value = 934.4286
$ws.Range("I116").Value = 963.1053000000001
$ws.Range("J116").Value = 662
$ws.Range("K116").Value = 963.1053000000001
$ws.Range("L116").Value = 662
$ws.Range("M116").Value = 1330.8947
$ws.Range("N116").Value = -5250
$ws.Range("H122").Value = 3507.6667
$ws.Range("I122").Value = 3507.6667
$ws.Range("K122").Value = 10523.0001
$ws.Range("M122").Value = -8073.000100000001
$ws.Range("H132").Value = 4261.6313
$ws.Range("I132").Value = 2822.818
$ws.Range("J132").Value = 6240
$ws.Range("K132").Value = 8468.454000000002
$ws.Range("L132").Value = 18720
$ws.Range("M132").Value = -5938.454000000002
$ws.Range("N132").Value = -23780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 934.4286
$ws.Range("I3").Value = 963.1053000000001
$ws.Range("J3").Value = 662
$ws.Range("K3").Value = 963.1053000000001
$ws.Range("L3").Value = 662
$ws.Range("M3").Value = -849.1053000000001
$ws.Range("N3").Value = -890

$ws.Range("H33").Value = 4250
$ws.Range("I33").Value = 4250
$ws.Range("K33").Value = 4250
$ws.Range("M33").Value = -3914
$ws.Range("H75").Value = 49999
$ws.Range("I75").Value = 49999
$ws.Range("K75").Value = 49999
$ws.Range("M75").Value = -49063
$ws.Range("H78").Value = 49999
$ws.Range("I78").Value = 49999
$ws.Range("K78").Value = 149997
$ws.Range("M78").Value = -145317
$ws.Range("H80").Value = 32161.385
$ws.Range("I80").Value = 38387.5
$ws.Range("J80").Value = 22199.6
$ws.Range("K80").Value = 38387.5
$ws.Range("L80").Value = 22199.6
$ws.Range("M80").Value = -37389.5
$ws.Range("N80").Value = -24195.6
$ws.Range("H82").Value = 52031.848
$ws.Range("J82").Value = 96660.664
$ws.Range("L82").Value = 96660.664
$ws.Range("N82").Value = -97426.664
$ws.Range("H83").Value = 32161.385
$ws.Range("I83").Value = 38387.5
$ws.Range("J83").Value = 22199.6
$ws.Range("K83").Value = 191937.5
$ws.Range("L83").Value = 110998
$ws.Range("M83").Value = -186945.5
$ws.Range("N83").Value = -120982
$ws.Range("H85").Value = 52031.848
$ws.Range("J85").Value = 96660.664
$ws.Range("L85").Value = 96660.664
$ws.Range("N85").Value = -99312.664
$ws.Range("H86").Value = 1849.8948
$ws.Range("I86").Value = 1867.8572
$ws.Range("J86").Value = 1799.6
$ws.Range("K86").Value = 1867.8572
$ws.Range("L86").Value = 1799.6
$ws.Range("M86").Value = -744.8571999999999
$ws.Range("N86").Value = -4045.6
$ws.Range("H89").Value = 1849.8948
$ws.Range("I89").Value = 1867.8572
$ws.Range("J89").Value = 1799.6
$ws.Range("K89").Value = 9339.286
$ws.Range("L89").Value = 8998
$ws.Range("M89").Value = -3723.286
$ws.Range("N89").Value = -20230
$ws.Range("H94").Value = 960.5294
$ws.Range("I94").Value = 853.5161000000001
$ws.Range("J94").Value = 2066.3333
$ws.Range("K94").Value = 853.5161000000001
$ws.Range("L94").Value = 2066.3333
$ws.Range("M94").Value = -402.5161000000001
$ws.Range("N94").Value = -2968.3333
$ws.Range("H97").Value = 8966.6
$ws.Range("I97").Value = 8715.625
$ws.Range("J97").Value = 9970.5
$ws.Range("K97").Value = 8715.625
$ws.Range("L97").Value = 9970.5
$ws.Range("M97").Value = -7724.625
$ws.Range("N97").Value = -11952.5
$ws.Range("H103").Value = 35096.8
$ws.Range("J103").Value = 35096.8
$ws.Range("L103").Value = 35096.8
$ws.Range("N103").Value = -37440.8
$ws.Range("H105").Value = 3490.4092
$ws.Range("I105").Value = 3458.111
$ws.Range("K105").Value = 3458.111
$ws.Range("M105").Value = -1711.111
$ws.Range("H107").Value = 2741.182
$ws.Range("I107").Value = 2849.9
$ws.Range("J107").Value = 1654
$ws.Range("K107").Value = 2849.9
$ws.Range("L107").Value = 1654
$ws.Range("M107").Value = -929.9000000000001
$ws.Range("N107").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 766.75
$ws.Range("I16").Value = 726
$ws.Range("J16").Value = 823.8
$ws.Range("K16").Value = 726
$ws.Range("L16").Value = 823.8
$ws.Range("M16").Value = -439
$ws.Range("N16").Value = -1397.8
$ws.Range("H29").Value = 12666.667
$ws.Range("J29").Value = 12666.667
$ws.Range("L29").Value = 12666.667
$ws.Range("N29").Value = -13252.667
$ws.Range("H58").Value = 3435.8667
$ws.Range("I58").Value = 3090
$ws.Range("K58").Value = 3090
$ws.Range("M58").Value = -2887
$ws.Range("H74").Value = 99749.5
$ws.Range("J74").Value = 99667
$ws.Range("L74").Value = 99667
$ws.Range("N74").Value = -101415
$ws.Range("H77").Value = 99749.5
$ws.Range("J77").Value = 99667
$ws.Range("L77").Value = 299001
$ws.Range("N77").Value = -307737
$ws.Range("H99").Value = 2537
$ws.Range("I99").Value = 2489.6667
$ws.Range("J99").Value = 2750
$ws.Range("K99").Value = 2489.6667
$ws.Range("L99").Value = 2750
$ws.Range("M99").Value = -991.6667000000002
$ws.Range("N99").Value = -5746
$ws.Range("H107").Value = 2316.2727
$ws.Range("I107").Value = 2113.5
$ws.Range("K107").Value = 2113.5
$ws.Range("M107").Value = -193.5
$ws.Range("H113").Value = 766.75
$ws.Range("I113").Value = 726
$ws.Range("J113").Value = 823.8
$ws.Range("K113").Value = 726
$ws.Range("L113").Value = 823.8
$ws.Range("M113").Value = 1444
$ws.Range("N113").Value = -5163.8
$ws.Range("H122").Value = 3751.0356
$ws.Range("I122").Value = 3559.842
$ws.Range("J122").Value = 4154.6665
$ws.Range("K122").Value = 10679.526
$ws.Range("L122").Value = 12463.9995
$ws.Range("M122").Value = -8229.526
$ws.Range("N122").Value = -17363.9995
$ws.Range("H126").Value = 2537
$ws.Range("I126").Value = 2489.6667
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 7469.000100000001
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -4999.000100000001
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 5511.263
$ws.Range("I132").Value = 3264.4285
$ws.Range("K132").Value = 9793.2855
$ws.Range("M132").Value = -7263.2855
$ws.Range("H134").Value = 14708452
$ws.Range("I134").Value = 21741134
$ws.Range("J134").Value = 3752.9092
$ws.Range("K134").Value = 65223402
$ws.Range("L134").Value = 11258.7276
$ws.Range("M134").Value = -65220867
$ws.Range("N134").Value = -16328.7276
$ws.Range("H136").Value = 3435.8667
$ws.Range("I136").Value = 3090
$ws.Range("K136").Value = 9270
$ws.Range("M136").Value = -6720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2026.375
$ws.Range("I5").Value = 2441.2
$ws.Range("J5").Value = 1335
$ws.Range("K5").Value = 7323.599999999999
$ws.Range("L5").Value = 4005
$ws.Range("M5").Value = -7211.599999999999
$ws.Range("N5").Value = -4229
$ws.Range("H55").Value = 7133.1665
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 9949.75
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 29849.25
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -30203.25
$ws.Range("H59").Value = 3750
$ws.Range("I59").Value = 3750
$ws.Range("K59").Value = 11250
$ws.Range("M59").Value = -10710
$ws.Range("H68").Value = 797.2727
$ws.Range("I68").Value = 610.25
$ws.Range("J68").Value = 1296
$ws.Range("K68").Value = 1830.75
$ws.Range("L68").Value = 3888
$ws.Range("M68").Value = -1019.75
$ws.Range("N68").Value = -5510
$ws.Range("H71").Value = 797.2727
$ws.Range("I71").Value = 610.25
$ws.Range("J71").Value = 1296
$ws.Range("K71").Value = 5492.25
$ws.Range("L71").Value = 11664
$ws.Range("M71").Value = -1436.25
$ws.Range("N71").Value = -19776
$ws.Range("H92").Value = 892.1667
$ws.Range("I92").Value = 1056
$ws.Range("K92").Value = 3168
$ws.Range("M92").Value = -1920
$ws.Range("H135").Value = 2026.375
$ws.Range("I135").Value = 2441.2
$ws.Range("J135").Value = 1335
$ws.Range("K135").Value = 21970.8
$ws.Range("L135").Value = 12015
$ws.Range("M135").Value = -19435.8
$ws.Range("N135").Value = -17085
$ws.Range("H140").Value = 1819.3103
$ws.Range("I140").Value = 1428.2941
$ws.Range("J140").Value = 2373.25
$ws.Range("K140").Value = 4284.8823
$ws.Range("L140").Value = 7119.75
$ws.Range("M140").Value = 895.1176999999998
$ws.Range("N140").Value = -17479.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 586.4761999999999
$ws.Range("I2").Value = 763.8125
$ws.Range("K2").Value = 763.8125
$ws.Range("M2").Value = -650.8125
$ws.Range("H26").Value = 32800
$ws.Range("J26").Value = 32800
$ws.Range("L26").Value = 32800
$ws.Range("N26").Value = -33360
$ws.Range("H50").Value = 32800
$ws.Range("J50").Value = 32800
$ws.Range("L50").Value = 32800
$ws.Range("N50").Value = -33796
$ws.Range("H70").Value = 8727.77
$ws.Range("I70").Value = 7932
$ws.Range("J70").Value = 8966.5
$ws.Range("K70").Value = 7932
$ws.Range("L70").Value = 8966.5
$ws.Range("M70").Value = -7662
$ws.Range("N70").Value = -9506.5
$ws.Range("H73").Value = 8727.77
$ws.Range("I73").Value = 7932
$ws.Range("J73").Value = 8966.5
$ws.Range("K73").Value = 7932
$ws.Range("L73").Value = 8966.5
$ws.Range("M73").Value = -6996
$ws.Range("N73").Value = -10838.5
$ws.Range("H97").Value = 579.3333
$ws.Range("I97").Value = 465.1111
$ws.Range("J97").Value = 807.7778
$ws.Range("K97").Value = 465.1111
$ws.Range("L97").Value = 807.7778
$ws.Range("M97").Value = 30.88889999999998
$ws.Range("N97").Value = -1799.7778
$ws.Range("H122").Value = 3353.8
$ws.Range("I122").Value = 1562.0769
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 4686.2307
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2236.2307
$ws.Range("N122").Value = -49900
$ws.Range("H132").Value = 3615.3809
$ws.Range("I132").Value = 1938.5333
$ws.Range("J132").Value = 7807.5
$ws.Range("K132").Value = 5815.5999
$ws.Range("L132").Value = 23422.5
$ws.Range("M132").Value = -3285.5999
$ws.Range("N132").Value = -28482.5
$ws.Range("H134").Value = 45694.777
$ws.Range("J134").Value = 45694.777
$ws.Range("L134").Value = 137084.331
$ws.Range("N134").Value = -142154.331
$ws.Range("H136").Value = 8763
$ws.Range("J136").Value = 8763
$ws.Range("L136").Value = 26289
$ws.Range("N136").Value = -31389

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1485.9584
$ws.Range("I22").Value = 836.7857
$ws.Range("J22").Value = 2394.8
$ws.Range("K22").Value = 836.7857
$ws.Range("L22").Value = 2394.8
$ws.Range("M22").Value = -541.7857
$ws.Range("N22").Value = -2984.8
$ws.Range("H27").Value = 1485.9584
$ws.Range("I27").Value = 836.7857
$ws.Range("J27").Value = 2394.8
$ws.Range("K27").Value = 836.7857
$ws.Range("L27").Value = 2394.8
$ws.Range("M27").Value = -729.7857
$ws.Range("N27").Value = -2608.8
$ws.Range("H46").Value = 2833.8948
$ws.Range("I46").Value = 908.25
$ws.Range("K46").Value = 908.25
$ws.Range("M46").Value = -720.25
$ws.Range("H61").Value = 2387.5715
$ws.Range("I61").Value = 2462.3157
$ws.Range("K61").Value = 2462.3157
$ws.Range("M61").Value = -2260.3157
$ws.Range("H93").Value = 1008.75
$ws.Range("I93").Value = 776.5625
$ws.Range("K93").Value = 776.5625
$ws.Range("M93").Value = 471.4375
$ws.Range("H113").Value = 2387.5715
$ws.Range("I113").Value = 2462.3157
$ws.Range("K113").Value = 2462.3157
$ws.Range("M113").Value = -292.3157000000001
$ws.Range("H122").Value = 5809.4165
$ws.Range("J122").Value = 7202
$ws.Range("L122").Value = 21606
$ws.Range("N122").Value = -26506
$ws.Range("H132").Value = 2226.4856
$ws.Range("I132").Value = 1388.5
$ws.Range("J132").Value = 3644.6155
$ws.Range("K132").Value = 4165.5
$ws.Range("L132").Value = 10933.8465
$ws.Range("M132").Value = -1635.5
$ws.Range("N132").Value = -15993.8465
$ws.Range("H136").Value = 3157.2
$ws.Range("I136").Value = 2363.4167
$ws.Range("J136").Value = 6332.3335
$ws.Range("K136").Value = 7090.250100000001
$ws.Range("L136").Value = 18997.0005
$ws.Range("M136").Value = -4540.250100000001
$ws.Range("N136").Value = -24097.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 14875
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 14500
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 14500
$ws.Range("M5").Value = -14888
$ws.Range("N5").Value = -14724
$ws.Range("H62").Value = 7703.3076
$ws.Range("I62").Value = 6774.1665
$ws.Range("J62").Value = 8499.714
$ws.Range("K62").Value = 6774.1665
$ws.Range("L62").Value = 8499.714
$ws.Range("M62").Value = -6150.1665
$ws.Range("N62").Value = -9747.714
$ws.Range("H65").Value = 7703.3076
$ws.Range("I65").Value = 6774.1665
$ws.Range("J65").Value = 8499.714
$ws.Range("K65").Value = 33870.8325
$ws.Range("L65").Value = 42498.57
$ws.Range("M65").Value = -30750.8325
$ws.Range("N65").Value = -48738.57
$ws.Range("H97").Value = 75518
$ws.Range("J97").Value = 75518
$ws.Range("L97").Value = 75518
$ws.Range("N97").Value = -77500
$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178
$ws.Range("H122").Value = 2422
$ws.Range("I122").Value = 2633.8635
$ws.Range("J122").Value = 1756.1428
$ws.Range("K122").Value = 7901.5905
$ws.Range("L122").Value = 5268.428400000001
$ws.Range("M122").Value = -5451.5905
$ws.Range("N122").Value = -10168.4284
$ws.Range("H132").Value = 3326.7878
$ws.Range("I132").Value = 2302.9312
$ws.Range("K132").Value = 6908.7936
$ws.Range("M132").Value = -4378.7936
$ws.Range("H136").Value = 2960.3704
$ws.Range("I136").Value = 2478.762
$ws.Range("K136").Value = 7436.286
$ws.Range("M136").Value = -4886.286
